$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 - this shifts the existing rows 42-44
# (the three most-recent "Arica y Parinacota" / "Perú" price records) down
# to rows 43-45, preserving all their original values and formatting.
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42 with the new weekly price record.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44746
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100108
$ws.Range("H42").Value = "Tropicales y subtropicales"
$ws.Range("I42").Value = 100108003
$ws.Range("J42").Value = "Maracuyá"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 70
$ws.Range("N42").Value = 30000
$ws.Range("O42").Value = 34000
$ws.Range("P42").Value = 31714
$ws.Range("Q42").Value = "$/caja 18 kilos"
$ws.Range("R42").Value = "Región de Arica y Parinacota"
$ws.Range("S42").Value = 1762
$ws.Range("T42").Value = 18
